$d = $word.ActiveDocument

$replacements = @(
    @{old = "42×79=3318"; new = "15×60=900"},
    @{old = "77×60=4620"; new = "25×38=950"},
    @{old = "58×71=4118"; new = "73×45=3285"},
    @{old = "45×42=1890"; new = "66×54=3564"},
    @{old = "45×26=1170"; new = "53×17=901"},
    @{old = "29×14=406";  new = "49×69=3381"},
    @{old = "47×86=4042"; new = "84×76=6384"},
    @{old = "46×35=1610"; new = "75×14=1050"},
    @{old = "76×45=3420"; new = "78×43=3354"},
    @{old = "52×16=832";  new = "64×91=5824"},
    @{old = "99×31=3069"; new = "25×67=1675"},
    @{old = "89×16=1424"; new = "33×76=2508"},
    @{old = "27×23=621";  new = "45×95=4275"},
    @{old = "55×49=2695"; new = "32×28=896"},
    @{old = "72×75=5400"; new = "90×74=6660"},
    @{old = "48×16=768";  new = "86×29=2494"},
    @{old = "60×19=1140"; new = "52×60=3120"},
    @{old = "44×58=2552"; new = "43×51=2193"},
    @{old = "50×20=1000"; new = "94×64=6016"},
    @{old = "78×15=1170"; new = "82×13=1066"},
    @{old = "83×34=2822"; new = "30×50=1500"},
    @{old = "59×66=3894"; new = "62×39=2418"},
    @{old = "72×47=3384"; new = "58×36=2088"},
    @{old = "66×53=3498"; new = "81×88=7128"},
    @{old = "28×89=2492"; new = "15×93=1395"}
)

foreach ($rep in $replacements) {
    $d.Content.Find.Execute($rep.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $rep.new, 2)
}
